$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = 0; B = 'Feldúlva találták a sírhelyeket'; C = 'Feldúlva találták a sírhelyeket, illetve azok környékét a Vasláb községhez tartozó hevederi temetőben. A nyomok alapján medvejárásra gyanakodnak.'; D = '[''Vasláb'', ''Heveder'']'; E = '[''Vasláb'', ''Heveder'']' },
    @{ A = 1; B = 'Elutasította medveügyben a Zetelaki Területtulajdonosi Társulás keresetét a táblabíróság'; C = 'Elutasította a Marosvásárhelyi Táblabíróság a Zetelaki Területtulajdonosi Társulás keresetét, amelyet a társulás a Környezetvédelmi Minisztérium ellen indított a vadgazdálkodási szabályozások alkalmazásának elmulasztása miatt. A társulás vezetője fellebbezést tervez.'; D = '[''Marosvásárhely'', ''Zetelaka'']'; E = '[''Marosvásárhely'', ''Zetelaka'']' },
    @{ A = 2; B = 'Ki akadályozta meg, hogy a barnamedve lekerüljön a szigorúan védett állatfajok listájáról Romániában?'; C = 'Az állatvédők és az Európai Bizottság akadályozta meg, hogy a barnamedve öt évre lekerüljön a szigorúan védett állatfajok listájáról Romániában – állítja Benkő Erika RMDSZ-képviselő. '; D = $null; E = $null },
    @{ A = 3; B = 'Pénzt ígér a miniszter a medvék által veszélyeztetett települések védelmére'; C = 'A háromszéki Zabolán tartott terepszemlét Costel Alexe környezetvédelmi miniszter, aki a látogatást követően arról számol be, hogy körvonalazódott egy olyan finanszírozási program, amely lehetővé teszi a medvék által veszélyeztetett települések védelmét.'; D = 'Zabola'; E = 'Zabola' },
    @{ A = 4; B = 'Hargita megye: megvan az év első 112-s medveészlelése'; C = 'Medvét látott a kertjében egy parajdi férfi a Sóhát utcában szombaton délután. A nagyvadat a gyümölcsfáknál fedezte fel, mintegy száz méterre a lakóháztól. Az esetet a 112-n jelentette, a helyszínre egy csendőri és egy mentőegység szállt ki.'; D = 'Parajd'; E = 'Parajd' },
    @{ A = 5; B = 'Aktívak a medvék Székelyudvarhely környékén'; C = 'Noha a magasabban fekvő térségekben már téli álmot alszanak a medvék, más területeken ez nem így van. Székelyudvarhely környékén például legalább tizenegy medve aktív jelenleg is, ezért a vadászok óvatosságra intenek.'; D = 'Székelyudvarhely'; E = 'Székelyudvarhely' },
    @{ A = 6; B = 'Lemondott a vadásztársaság az emberre támadó hidegkúti medve kilövéséről'; C = 'Megúszta a kilövést a Hidegkúton emberre támadó medve, az illetékes vadásztársaságnál lemondtak arról, hogy a vad ártalmatlanítására rendkívüli jóváhagyást igényeljenek a környezetvédelmi minisztériumtól.'; D = 'Hidegkút'; E = 'Hidegkút' },
    @{ A = 7; B = 'Amíg elérhető közelségben van az ételmaradék, addig a medvék jelenlétére is számítani kell'; C = 'A szeméttárolók vonzzák a székelyudvarhelyi Cserehát lakónegyedbe az aktív nagyvadakat, ezért a Nagy-Küküllő Vadász- és Sporthorgász Egyesület medvebiztos kukákat rendelt, amelyeket a szemételszállító vállalattal egyeztetve helyezne ki. A medvék befogásával is próbálkoznak.'; D = 'Székelyudvarhely'; E = 'Székelyudvarhely' },
    @{ A = 8; B = 'Medveradar: Zetelakán és Farkaslakán voltak a legaktívabbak tavaly a nagyvadak'; C = 'A kezdeti fellendülés után kissé lankadt az aktivitás a székelyföldi medveradar és -térkép néven emlegetett medveészlelő portálon, amelyet tavaly márciusban hozott létre Csala Dénes adatblogger. Ettől eltekintve a több mint 700 bejegyzést számláló medvetérkép a legszerteágazóbb adatbázisnak számít.'; D = '[''Haraly'', ''Zetelaka'', ''Farkaslaka'']'; E = '[''Haraly'', ''Zetelaka'', ''Farkaslaka'']' },
    @{ A = 9; B = 'Gyergyószentmiklós utcáin kóborolt egy medve'; C = 'Egy városszéli üzemanyagtöltő állomásnál felbukkant medve miatt riasztották keddre virradóan a hatóságokat Gyergyószentmiklóson. Míg a csendőrök a nagyvadat kutatták, a vészhelyzeti felügyelőség a Ro-Alert rendszeren figyelmeztette a lakókat a veszélyre. A medvét végül megtalálták és elűzték.'; D = 'Gyergyószentmiklós'; E = 'Gyergyószentmiklós' }
)

# Row 11 is a brand-new row; give its A cell the same number formatting/style
# (bold, bordered, centered) that the other index cells in column A already use,
# by copying the format from A10 before filling in the values.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    if ($r.D -ne $null) {
        $ws.Cells.Item($rowNum, 4).Value = $r.D
    } else {
        $ws.Cells.Item($rowNum, 4).ClearContents()
    }
    if ($r.E -ne $null) {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    } else {
        $ws.Cells.Item($rowNum, 5).ClearContents()
    }
}
